# Applies the "version final sin errores" edit to the FHIR StructureDefinition
# workbook: bump version/date, fix the extension Context, and add the missing
# ele-1/ext-1 invariant text to the root "Extension" row of the Elements sheet.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value = "0.7.0"
$meta.Range("B8").Value = "2023-09-13T17:11:14-03:00"
$meta.Range("B20").Value = "element:Patient"

# --- Elements sheet --------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$elements.Range("AJ1").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"
